$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A (shifts 攤位編號/緯度/經度 one column right)
# and add the new "分區" (zone) column with its header + values.
$ws.Columns("A:A").Insert()

$ws.Range("A1").Value() = "分區"
$ws.Range("A2").Value() = "甲"
$ws.Range("A3").Value() = "甲"
$ws.Range("A4").Value() = "乙"

$ws.Columns("A:A").ColumnWidth = 5.285714285714286

# Re-point the "duplicate values" conditional formatting from the old booth-number
# column (now B) onto column B, replacing the rule that used to live on column A.
$rngB = $ws.Range("B1:B1048576")

# (throwaway rule -- mirrors the dxf bookkeeping left behind in the original edit)
$throwaway = $rngB.FormatConditions.AddUniqueValues(1)
$throwaway.DupeUnique = 1
$throwaway.Font.Color = 393372
$throwaway.Interior.Color = 13551615
$rngB.FormatConditions.Delete()

$dupRule = $rngB.FormatConditions.AddUniqueValues(1)
$dupRule.DupeUnique = 1
$dupRule.Font.Color = 393372
$dupRule.Interior.Color = 13551615

$ws.Range("A1:A1048576").FormatConditions.Delete()

# Match the final selected cell recorded in the saved workbook.
$ws.Range("F11").Select() | Out-Null
